# Golf.xlsx fixture update: "Add form_txt check for editable fields"
# Adds a new "Fairway" (hit Y/N) row and a new "80 <" (score bucket) row to
# the most recent round on the "Knight's play 10-18" sheet, and appends a
# brand-new round (with Fairway/80< rows too) to the "Lochmere" sheet.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Sheet: "Knight's play 10-18" (2nd tab)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Round date (existing blank cell A45 already carries the date style).
$ws2.Range("A45").Value = 44402

# Black / White / Blue header cell + per-hole scores for this round.
$ws2.Range("A46").Value = "Black"
$ws2.Range("D46").Value = 3
$ws2.Range("E46").Value = 5
$ws2.Range("F46").Value = 4
$ws2.Range("G46").Value = 4
$ws2.Range("H46").Value = 4
$ws2.Range("I46").Value = 4
$ws2.Range("J46").Value = 4
$ws2.Range("K46").Value = 3
$ws2.Range("L46").Value = 8
$ws2.Range("M46").Formula = "=SUM(D46:L46)"

# Putts row.
$ws2.Range("A47").Value = "putts"
$ws2.Range("D47").Value = 1
$ws2.Range("E47").Value = 2
$ws2.Range("F47").Value = 2
$ws2.Range("G47").Value = 2
$ws2.Range("H47").Value = 2
$ws2.Range("I47").Value = 2
$ws2.Range("J47").Value = 2
$ws2.Range("K47").Value = 2
$ws2.Range("L47").Value = 3
$ws2.Range("M47").Formula = "=SUM(D47:L47)"

# Penalties row.
$ws2.Range("A48").Value = "penalties"
$ws2.Range("E48").Value = "W"
$ws2.Range("K48").Value = "M"
$ws2.Range("L48").Value = "W"

# New "Fairway" row (hit in regulation?): copy the bold label style from the
# existing "Black" label cell, then fill in the X markers.
$ws2.Range("A46").Copy()
$ws2.Range("A49").PasteSpecial($xlPasteFormats)
$ws2.Range("A49").Value = "Fairway"
foreach ($col in @("D","E","F","G","H","I","J","K","L")) {
    $ws2.Range($col + "49").Value = "X"
}

# New "80 <" row (score bucket counts).
$ws2.Range("A46").Copy()
$ws2.Range("A50").PasteSpecial($xlPasteFormats)
$ws2.Range("A50").Value = "80 <"
$ws2.Range("D50").Value = 2
$ws2.Range("E50").Value = 4
$ws2.Range("F50").Value = 3
$ws2.Range("G50").Value = 4
$ws2.Range("H50").Value = 3
$ws2.Range("I50").Value = 3
$ws2.Range("J50").Value = 3
$ws2.Range("K50").Value = 2
$ws2.Range("L50").Value = 5

$ws2.Range("A49:M50").Select()

# ---------------------------------------------------------------------------
# Sheet: "Lochmere" (4th tab) - append a brand-new round.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Date row, copying the date-style cell from the previous round (A40).
$ws4.Range("A40").Copy()
$ws4.Range("A45").PasteSpecial($xlPasteFormats)
$ws4.Range("A45").Value = 44397

# Tee-time cells: apply the (new) time number format, then set the values.
$ws4.Range("B45").NumberFormat = "h:mm AM/PM"
$ws4.Range("B45").Value = 0.45833333333333331
$ws4.Range("C45").NumberFormat = "h:mm AM/PM"
$ws4.Range("C45").Value = 0.63541666666666663

# Blue tee header + per-hole scores (front nine + back nine), copying styles
# from the previous round's equivalent row (41).
$ws4.Range("A41").Copy()
$ws4.Range("A46").PasteSpecial($xlPasteFormats)
$ws4.Range("A46").Value = "Blue"
$ws4.Range("B41").Copy()
$ws4.Range("B46").PasteSpecial($xlPasteFormats)
$ws4.Range("C41").Copy()
$ws4.Range("C46").PasteSpecial($xlPasteFormats)

$ws4.Range("D46").Value = 5
$ws4.Range("E46").Value = 6
$ws4.Range("F46").Value = 4
$ws4.Range("G46").Value = 8
$ws4.Range("H46").Value = 5
$ws4.Range("I46").Value = 5
$ws4.Range("J46").Value = 6
$ws4.Range("K41").Copy()
$ws4.Range("K46").PasteSpecial($xlPasteFormats)
$ws4.Range("K46").Value = 7
$ws4.Range("L41").Copy()
$ws4.Range("L46").PasteSpecial($xlPasteFormats)
$ws4.Range("L46").Value = 8
$ws4.Range("M41").Copy()
$ws4.Range("M46").PasteSpecial($xlPasteFormats)
$ws4.Range("M46").Formula = "=SUM(D46:L46)"

$ws4.Range("N46").Value = 8
$ws4.Range("O46").Value = 5
$ws4.Range("P46").Value = 10
$ws4.Range("Q46").Value = 7
$ws4.Range("R46").Value = 6
$ws4.Range("S46").Value = 5
$ws4.Range("T46").Value = 7
$ws4.Range("U41").Copy()
$ws4.Range("U46").PasteSpecial($xlPasteFormats)
$ws4.Range("U46").Value = 4
$ws4.Range("V41").Copy()
$ws4.Range("V46").PasteSpecial($xlPasteFormats)
$ws4.Range("V46").Value = 5
$ws4.Range("W41").Copy()
$ws4.Range("W46").PasteSpecial($xlPasteFormats)
$ws4.Range("W46").Formula = "=SUM(N46:V46)"
$ws4.Range("X41").Copy()
$ws4.Range("X46").PasteSpecial($xlPasteFormats)
$ws4.Range("X46").Formula = "=W46+M46"

# Putts row, copying styles from row 42.
$ws4.Range("A42").Copy()
$ws4.Range("A47").PasteSpecial($xlPasteFormats)
$ws4.Range("A47").Value = "putts"

$ws4.Range("D47").Value = 2
$ws4.Range("E47").Value = 3
$ws4.Range("F47").Value = 1
$ws4.Range("G47").Value = 2
$ws4.Range("H47").Value = 1
$ws4.Range("I47").Value = 2
$ws4.Range("J47").Value = 3
$ws4.Range("K42").Copy()
$ws4.Range("K47").PasteSpecial($xlPasteFormats)
$ws4.Range("K47").Value = 2
$ws4.Range("L42").Copy()
$ws4.Range("L47").PasteSpecial($xlPasteFormats)
$ws4.Range("L47").Value = 3
$ws4.Range("M42").Copy()
$ws4.Range("M47").PasteSpecial($xlPasteFormats)
$ws4.Range("M47").Formula = "=SUM(D47:L47)"

$ws4.Range("N47").Value = 3
$ws4.Range("O47").Value = 1
$ws4.Range("P47").Value = 2
$ws4.Range("Q47").Value = 2
$ws4.Range("R47").Value = 2
$ws4.Range("S47").Value = 2
$ws4.Range("T47").Value = 3
$ws4.Range("U42").Copy()
$ws4.Range("U47").PasteSpecial($xlPasteFormats)
$ws4.Range("U47").Value = 2
$ws4.Range("V42").Copy()
$ws4.Range("V47").PasteSpecial($xlPasteFormats)
$ws4.Range("V47").Value = 2
$ws4.Range("W42").Copy()
$ws4.Range("W47").PasteSpecial($xlPasteFormats)
$ws4.Range("W47").Formula = "=SUM(N47:V47)"
$ws4.Range("X42").Copy()
$ws4.Range("X47").PasteSpecial($xlPasteFormats)
$ws4.Range("X47").Formula = "=W47+M47"

# Penalties row, copying styles from row 43.
$ws4.Range("A43").Copy()
$ws4.Range("A48").PasteSpecial($xlPasteFormats)
$ws4.Range("A48").Value = "penalties"

$ws4.Range("D48").Value = "M"
$ws4.Range("F48").Value = "W"
$ws4.Range("K43").Copy()
$ws4.Range("K48").PasteSpecial($xlPasteFormats)
$ws4.Range("O48").Value = "W"
$ws4.Range("P48").Value = " "
$ws4.Range("Q48").Value = "L"
$ws4.Range("R48").Value = "WW"
$ws4.Range("T48").Value = "W"
$ws4.Range("U43").Copy()
$ws4.Range("U48").PasteSpecial($xlPasteFormats)
$ws4.Range("U48").Value = "M"
$ws4.Range("V43").Copy()
$ws4.Range("V48").PasteSpecial($xlPasteFormats)

# New "Fairway" row (hit in regulation?).
$ws4.Range("A48").Copy()
$ws4.Range("A49").PasteSpecial($xlPasteFormats)
$ws4.Range("A49").Value = "Fairway"
foreach ($col in @("D","E","F","G","H","I","J","L")) {
    $ws4.Range($col + "49").Value = "N"
}
$ws4.Range("K49").Value = "Y"
foreach ($col in @("N","O","P","Q","T","U","V")) {
    $ws4.Range($col + "49").Value = "N"
}
$ws4.Range("R49").Value = "Y"
$ws4.Range("S49").Value = "Y"

# New "80 <" row (score bucket counts).
$ws4.Range("A48").Copy()
$ws4.Range("A50").PasteSpecial($xlPasteFormats)
$ws4.Range("A50").Value = "80 <"
$ws4.Range("D50").Value = 3
$ws4.Range("E50").Value = 4
$ws4.Range("F50").Value = 3
$ws4.Range("G50").Value = 5
$ws4.Range("H50").Value = 2
$ws4.Range("I50").Value = 3
$ws4.Range("J50").Value = 4
$ws4.Range("K50").Value = 3
$ws4.Range("L50").Value = 4
$ws4.Range("N50").Value = 5
$ws4.Range("O50").Value = 3
$ws4.Range("P50").Value = 4
$ws4.Range("Q50").Value = 5
$ws4.Range("R50").Value = 3
$ws4.Range("S50").Value = 3
$ws4.Range("T50").Value = 4
$ws4.Range("U50").Value = 3
$ws4.Range("V50").Value = 3

$ws4.Range("W50").Select()

# Make "Lochmere" the active tab (mirrors the tabSelected move in the diff).
$ws4.Activate()
